# Random Forest algoritmasi icin Final model egitimi tamamlandi.
# Update the "Random Forest" results column (C) on Sheet1 with the
# final trained model's metrics / predictions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Scratch cell used to build exact TEXT() values below.
$scratch = $ws.Range("Z1000")

function Set-TextValue($range, $text) {
    $scratch.Formula = "=TEXT(" + $text + ",""0.0000"")"
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Model name ---
$ws.Range("C1").Value  = "Random Forest"

# --- Numeric-looking score values: these must stay stored as TEXT      ---
# --- (matching the rest of the workbook) rather than being converted  ---
# --- to real numbers by the usual Value auto-typing. We build the     ---
# --- text in a scratch cell with TEXT(), then copy/paste the computed ---
# --- value only, which keeps the cell's existing text style/format.   ---
Set-TextValue $ws.Range("C24") "0.8926"
Set-TextValue $ws.Range("C25") "0.8958"
Set-TextValue $ws.Range("C26") "0.8952"
Set-TextValue $ws.Range("C27") "0.9550"

$scratch.ClearContents()

# --- Model id / timestamp ---
$ws.Range("C2").Value  = "anxiety_model_20250510_1009"

# --- Prediction labels (plain text, no numeric coercion issue) ---
$ws.Range("C4").Value  = [char]0x2192 + "Anksiyete (84.31%)"
$ws.Range("C6").Value  = [char]0x2192 + "Normal (14.91%)"
$ws.Range("C8").Value  = [char]0x2192 + "Anksiyete (69.31%)"
$ws.Range("C10").Value = [char]0x2192 + "Normal (5.23%)"
$ws.Range("C12").Value = [char]0x2192 + "Normal (30.06%)"
$ws.Range("C14").Value = [char]0x2192 + "Normal (24.60%)"
$ws.Range("C16").Value = [char]0x2192 + "Normal (41.52%)"
$ws.Range("C18").Value = [char]0x2192 + "Normal (8.27%)"
$ws.Range("C20").Value = [char]0x2192 + "Normal (38.21%)"
$ws.Range("C22").Value = [char]0x2192 + "Normal (5.13%)"

# The author left the selection on D11 after finishing the edit.
$ws.Range("D11").Select()
